# Insert a new weekly price record for "Femacal de La Calera - Apio" at
# row 371, pushing the existing rows 371-396 down to 372-397.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(371).Insert()

$ws.Cells.Item(371, 1).Value  = 3
$ws.Cells.Item(371, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(371, 3).Value  = "Coquimbo"
$ws.Cells.Item(371, 4).Value  = 44714
$ws.Cells.Item(371, 5).Value  = 5
$ws.Cells.Item(371, 6).Value  = 100112017
$ws.Cells.Item(371, 7).Value  = "Apio"
$ws.Cells.Item(371, 8).Value  = "Americana (o)"
$ws.Cells.Item(371, 9).Value  = "Primera"
$ws.Cells.Item(371, 10).Value = 280
$ws.Cells.Item(371, 11).Value = 9000
$ws.Cells.Item(371, 12).Value = 9500
$ws.Cells.Item(371, 13).Value = 9286
$ws.Cells.Item(371, 14).Value = "$/docena de matas"
$ws.Cells.Item(371, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(371, 16).Value = 1548
$ws.Cells.Item(371, 17).Value = 6
$ws.Cells.Item(371, 18).Value = "Hortaliza"
